$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.486.04"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.098.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.89"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.30"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.14%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.092.97"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.29%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +6.11%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.615.27"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.351.13"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.11"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.097.24"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.93"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.45"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.26"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.12%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.96"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +8.56%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.67"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.81"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.74"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.68%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0851"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.46%  "

$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.40"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +7.82%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.88%  "

$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.03"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.32"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "439.27"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.73"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0369"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.863.16"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.50%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.268"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.90"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.69"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.20%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.17"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.21%  "
